# Append the latest daily quotations row (2025-10-23) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 49

# Column A holds the date as a serial number, formatted like the rows above it.
$ws.Range("A$newRow").Value = 45953
$ws.Range("A$newRow").NumberFormat = $ws.Range("A48").NumberFormat

# Columns B:E hold the quotation values as text (comma decimal separator).
$ws.Range("B$newRow").Value = "21,7178"
$ws.Range("C$newRow").Value = "15,6198"
$ws.Range("D$newRow").Value = "15,4273"
$ws.Range("E$newRow").Value = "15,4273"
